$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Status cells for rule1 and rule2 to be quoted: "拒保"
$ws.Range("C10").Value = '"拒保"'
$ws.Range("C11").Value = '"拒保"'

# Update the active selection to C11
$ws.Range("C11").Select()
